# Update column G ("K") values in Sheet1 for rows 2-51.
# These values were regenerated to reflect a "K" (strikeout) count
# instead of the previous "Strike#" metric.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 4
    6  = 2
    7  = 0
    8  = 1
    9  = 1
    10 = 1
    11 = 2
    12 = 1
    13 = 2
    14 = 1
    15 = 1
    16 = 1
    17 = 2
    18 = 0
    19 = 0
    20 = 2
    21 = 2
    22 = 2
    23 = 2
    24 = 3
    25 = 1
    26 = 0
    27 = 0
    28 = 2
    29 = 1
    30 = 2
    31 = 3
    32 = 0
    33 = 1
    34 = 1
    35 = 1
    36 = 0
    37 = 3
    38 = 5
    39 = 2
    40 = 0
    41 = 1
    42 = 2
    43 = 0
    44 = 4
    45 = 1
    46 = 2
    47 = 1
    48 = 2
    49 = 1
    50 = 2
    51 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
